# Plot Meziadin total return by year and sibling reg point forecast
# with 25/75th percentile error bars: add "p25" and "p75" columns to
# Table1, populate the data, and update the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Data for the two new columns -----------------------------------
# row => (p25, p75), rows correspond to runyear 2012..2025 (sheet rows 2..15)
$p25 = @(287613, 236161, 327560, 278110, 285478, 112488, 137562, 133096, 133189, 119280, 775247, 407566, 475099, 711287)
$p75 = @(651471, 605663, 633701, 597477, 567755, 309076, 339044, 339139, 331119, 309427, 1107526, 703180, 774011, 993255)

# --- Expand the existing table from A1:B15 to A1:D15 -----------------
$newRange = $ws.Range("A1:D15")
$tbl.Resize($newRange)

# Name the two new header cells (this also renames the table columns)
$ws.Range("C1").Value = "p25"
$ws.Range("D1").Value = "p75"

# --- Fill in the data rows, matching the integer number format used
#     by the "runyear" column (cell style s="2", numFmtId "0") --------
for ($i = 0; $i -lt 14; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $p25[$i]
    $ws.Cells.Item($r, 4).Value = $p75[$i]
}

$ws.Range("C2:C15").NumberFormat = "0"
$ws.Range("D2:D15").NumberFormat = "0"

# --- Update the active selection on the sheet -------------------------
$ws.Range("I1").Select()
